$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Authors" column (E) holds pseudo-CSV records such as
#   [First%Last%email%n,                First%Last%email%n, ...]
# separated by runs of spaces after each comma. The source data was
# regenerated with one additional space inserted after every comma in
# those separators. Re-apply the same fix to every row's Authors cell.
for ($row = 2; $row -le 26; $row++) {
    $cell = $ws.Cells.Item($row, 5)   # column E
    $value = $cell.Value2
    if ($value -ne $null -and $value -ne "") {
        $updated = $value -replace ',( +)', ', $1'
        if ($updated -ne $value) {
            $cell.Value2 = $updated
        }
    }
}
